# Add a "browser" column to the DATA sheet, between "execute" and "username".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Insert a new column before column C (shifts username/password/fname right).
$ws.Columns("C").Insert()

# Header
$ws.Range("C1").Value = "browser"

# Data rows: chrome/firefox for executed tests, N/A for skipped ones.
$ws.Range("C2").Value = "chrome"
$ws.Range("C3").Value = "N/A"
$ws.Range("C4").Value = "firefox"
$ws.Range("C5").Value = "N/A"
$ws.Range("C6").Value = "firefox"
$ws.Range("C7").Value = "N/A"

# Match the column width style used in the final workbook (same as column B).
$ws.Columns("C").ColumnWidth = 6.666666666666667

# Restore the active sheet's selection to match the end state.
$ws.Range("C7").Select()
